# Auto-generated Excel COM-interop script
# Updates Leve profit/price tracking values across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled market-data runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 764.2941
$ws.Range("I28").Value = 530.2308
$ws.Range("J28").Value = 1525
$ws.Range("K28").Value = 530.2308
$ws.Range("L28").Value = 1525
$ws.Range("M28").Value = -45.23080000000004

$ws.Range("H32").Value = 10372.625
$ws.Range("I32").Value = 7631.6665
$ws.Range("J32").Value = 12017.2
$ws.Range("K32").Value = 7631.6665
$ws.Range("L32").Value = 12017.2
$ws.Range("M32").Value = -7305.6665

$ws.Range("H33").Value = 16667859
$ws.Range("I33").Value = 25000676
$ws.Range("J33").Value = 2224
$ws.Range("K33").Value = 25000676
$ws.Range("L33").Value = 2224
$ws.Range("M33").Value = -25000447

$ws.Range("H76").Value = 5081.125
$ws.Range("I76").Value = 5622.25
$ws.Range("J76").Value = 4540
$ws.Range("K76").Value = 5622.25
$ws.Range("L76").Value = 4540
$ws.Range("M76").Value = -5307.25
$ws.Range("N76").Value = -5170

$ws.Range("H79").Value = 5081.125
$ws.Range("I79").Value = 5622.25
$ws.Range("J79").Value = 4540
$ws.Range("K79").Value = 5622.25
$ws.Range("L79").Value = 4540
$ws.Range("M79").Value = -4530.25
$ws.Range("N79").Value = -6724

$ws.Range("H125").Value = 5372
$ws.Range("I125").Value = 1715
$ws.Range("J125").Value = 20000
$ws.Range("K125").Value = 15435
$ws.Range("L125").Value = 180000
$ws.Range("M125").Value = -12975
$ws.Range("N125").Value = -184920

$ws.Range("H135").Value = 1081.1852
$ws.Range("I135").Value = 925.8261
$ws.Range("J135").Value = 1974.5
$ws.Range("K135").Value = 8332.4349
$ws.Range("L135").Value = 17770.5
$ws.Range("M135").Value = -5797.4349
$ws.Range("N135").Value = -22840.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 994.9091
$ws.Range("I45").Value = 777.7143
$ws.Range("J45").Value = 1375
$ws.Range("K45").Value = 777.7143
$ws.Range("L45").Value = 1375
$ws.Range("M45").Value = -400.7143

$ws.Range("H61").Value = 7937.255
$ws.Range("I61").Value = 2195.9614
$ws.Range("J61").Value = 13908.2
$ws.Range("K61").Value = 2195.9614
$ws.Range("L61").Value = 13908.2
$ws.Range("M61").Value = -1983.9614

$ws.Range("H63").Value = 3902.6667
$ws.Range("I63").Value = 3902.6667
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3902.6667
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -3216.6667
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 3902.6667
$ws.Range("I66").Value = 3902.6667
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 19513.3335
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -16081.3335
$ws.Range("N66").ClearContents()

$ws.Range("H102").Value = 11666.421
$ws.Range("I102").Value = 1166.7333
$ws.Range("J102").Value = 51040.25
$ws.Range("K102").Value = 1166.7333
$ws.Range("L102").Value = 51040.25
$ws.Range("M102").Value = 455.2666999999999
$ws.Range("N102").Value = -54284.25

$ws.Range("H122").Value = 2315.9092
$ws.Range("I122").Value = 1603.6487
$ws.Range("J122").Value = 6080.7144
$ws.Range("K122").Value = 4810.9461
$ws.Range("L122").Value = 18242.1432
$ws.Range("M122").Value = -2360.9461
$ws.Range("N122").Value = -23142.1432

$ws.Range("H136").Value = 7937.255
$ws.Range("I136").Value = 2195.9614
$ws.Range("J136").Value = 13908.2
$ws.Range("K136").Value = 6587.8842
$ws.Range("L136").Value = 41724.60000000001
$ws.Range("M136").Value = -4037.8842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 31810.818
$ws.Range("I20").Value = 21382.291
$ws.Range("J20").Value = 44325.05
$ws.Range("K20").Value = 21382.291
$ws.Range("L20").Value = 44325.05
$ws.Range("M20").Value = -21135.291
$ws.Range("N20").Value = -44819.05

$ws.Range("H86").Value = 2312.6333
$ws.Range("I86").Value = 2008.5714
$ws.Range("J86").Value = 3022.111
$ws.Range("K86").Value = 2008.5714
$ws.Range("L86").Value = 3022.111
$ws.Range("M86").Value = -885.5714
$ws.Range("N86").Value = -5268.111

$ws.Range("H89").Value = 2312.6333
$ws.Range("I89").Value = 2008.5714
$ws.Range("J89").Value = 3022.111
$ws.Range("K89").Value = 10042.857
$ws.Range("L89").Value = 15110.555
$ws.Range("M89").Value = -4426.857
$ws.Range("N89").Value = -26342.555

$ws.Range("H94").Value = 5560
$ws.Range("I94").Value = 700
$ws.Range("J94").Value = 25000
$ws.Range("K94").Value = 700
$ws.Range("L94").Value = 25000
$ws.Range("M94").Value = -249

$ws.Range("H99").Value = 1490.7
$ws.Range("I99").Value = 1283.1765
$ws.Range("J99").Value = 2666.6667
$ws.Range("K99").Value = 1283.1765
$ws.Range("L99").Value = 2666.6667
$ws.Range("M99").Value = 214.8235

$ws.Range("H103").Value = 23080.75
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 23080.75
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 23080.75
$ws.Range("N103").Value = -25424.75

$ws.Range("H105").Value = 3166.5
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 3166.5
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 3166.5
$ws.Range("N105").Value = -6660.5

$ws.Range("H106").Value = 55000
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 55000
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 55000
$ws.Range("N106").Value = -57524

$ws.Range("H134").Value = 14253.852
$ws.Range("I134").Value = 7774.294
$ws.Range("J134").Value = 25269.1
$ws.Range("K134").Value = 23322.882
$ws.Range("L134").Value = 75807.29999999999
$ws.Range("M134").Value = -20787.882
$ws.Range("N134").Value = -80877.29999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 7546.3213
$ws.Range("I99").Value = 2638.3845
$ws.Range("J99").Value = 11799.866
$ws.Range("K99").Value = 2638.3845
$ws.Range("L99").Value = 11799.866
$ws.Range("M99").Value = -1140.3845
$ws.Range("N99").Value = -14795.866

$ws.Range("H105").Value = 13938.75
$ws.Range("I105").Value = 50405
$ws.Range("J105").Value = 1783.3334
$ws.Range("K105").Value = 50405
$ws.Range("L105").Value = 1783.3334
$ws.Range("M105").Value = -48658
$ws.Range("N105").Value = -5277.3334

$ws.Range("H126").Value = 7546.3213
$ws.Range("I126").Value = 2638.3845
$ws.Range("J126").Value = 11799.866
$ws.Range("K126").Value = 7915.1535
$ws.Range("L126").Value = 35399.598
$ws.Range("M126").Value = -5445.1535
$ws.Range("N126").Value = -40339.598

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("M68:N68").ClearContents()

$ws.Range("H69").Value = 4501.4443
$ws.Range("I69").Value = 2583.1667
$ws.Range("J69").Value = 8338
$ws.Range("K69").Value = 7749.500100000001
$ws.Range("L69").Value = 25014
$ws.Range("M69").Value = -6938.500100000001
$ws.Range("N69").Value = -26636

$ws.Range("H70").Value = 13237.2
$ws.Range("I70").Value = 11374.667
$ws.Range("J70").Value = 30000
$ws.Range("K70").Value = 34124.001
$ws.Range("L70").Value = 90000
$ws.Range("M70").Value = -33809.001

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("M71:N71").ClearContents()

$ws.Range("H72").Value = 4501.4443
$ws.Range("I72").Value = 2583.1667
$ws.Range("J72").Value = 8338
$ws.Range("K72").Value = 23248.5003
$ws.Range("L72").Value = 75042
$ws.Range("M72").Value = -19192.5003
$ws.Range("N72").Value = -83154

$ws.Range("H73").Value = 13237.2
$ws.Range("I73").Value = 11374.667
$ws.Range("J73").Value = 30000
$ws.Range("K73").Value = 34124.001
$ws.Range("L73").Value = 90000
$ws.Range("M73").Value = -33032.001

$ws.Range("H122").Value = 13451579
$ws.Range("I122").Value = 23358916
$ws.Range("J122").Value = 3544242.2
$ws.Range("K122").Value = 210230244
$ws.Range("L122").Value = 31898179.8
$ws.Range("M122").Value = -210227794
$ws.Range("N122").Value = -31903079.8

$ws.Range("H126").Value = 9187.375
$ws.Range("I126").Value = 5699.8
$ws.Range("J126").Value = 15000
$ws.Range("K126").Value = 17099.4
$ws.Range("L126").Value = 45000
$ws.Range("M126").Value = -12159.4
$ws.Range("N126").Value = -54880

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 121.666664
$ws.Range("I2").Value = 129
$ws.Range("J2").Value = 103.85714
$ws.Range("K2").Value = 129
$ws.Range("L2").Value = 103.85714
$ws.Range("M2").Value = -16
$ws.Range("N2").Value = -329.85714

$ws.Range("H80").Value = 21760.889
$ws.Range("I80").Value = 32333
$ws.Range("J80").Value = 16474.834
$ws.Range("K80").Value = 32333
$ws.Range("L80").Value = 16474.834
$ws.Range("M80").Value = -31335
$ws.Range("N80").Value = -18470.834

$ws.Range("H83").Value = 21760.889
$ws.Range("I83").Value = 32333
$ws.Range("J83").Value = 16474.834
$ws.Range("K83").Value = 161665
$ws.Range("L83").Value = 82374.17
$ws.Range("M83").Value = -156673
$ws.Range("N83").Value = -92358.17

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 5493.857
$ws.Range("I68").Value = 1657.8572
$ws.Range("J68").Value = 9329.857
$ws.Range("K68").Value = 1657.8572
$ws.Range("L68").Value = 9329.857
$ws.Range("M68").Value = -908.8571999999999
$ws.Range("N68").Value = -10827.857

$ws.Range("H71").Value = 5493.857
$ws.Range("I71").Value = 1657.8572
$ws.Range("J71").Value = 9329.857
$ws.Range("K71").Value = 8289.286
$ws.Range("L71").Value = 46649.285
$ws.Range("M71").Value = -4545.286
$ws.Range("N71").Value = -54137.285

$ws.Range("H82").Value = 3847.2632
$ws.Range("I82").Value = 2318.625
$ws.Range("J82").Value = 4959
$ws.Range("K82").Value = 2318.625
$ws.Range("L82").Value = 4959
$ws.Range("M82").Value = -1957.625
$ws.Range("N82").Value = -5681

$ws.Range("H85").Value = 3847.2632
$ws.Range("I85").Value = 2318.625
$ws.Range("J85").Value = 4959
$ws.Range("K85").Value = 2318.625
$ws.Range("L85").Value = 4959
$ws.Range("M85").Value = -1070.625
$ws.Range("N85").Value = -7455

$ws.Range("H102").Value = 75000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 75000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 75000
$ws.Range("N102").Value = -81490

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2099.875
$ws.Range("I96").Value = 1759.8
$ws.Range("J96").Value = 2254.4546
$ws.Range("K96").Value = 1759.8
$ws.Range("L96").Value = 2254.4546
$ws.Range("M96").Value = -386.8
$ws.Range("N96").Value = -5000.4546

$ws.Range("H124").Value = 53036.332
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 53036.332
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 53036.332
$ws.Range("N124").Value = -62856.332

$ws.Range("H132").Value = 6338.515
$ws.Range("I132").Value = 1855.2174
$ws.Range("J132").Value = 16650.1
$ws.Range("K132").Value = 5565.6522
$ws.Range("L132").Value = 49950.3
$ws.Range("M132").Value = -3035.6522
$ws.Range("N132").Value = -55010.3

$ws.Range("H136").Value = 8731.459000000001
$ws.Range("I136").Value = 2642.5
$ws.Range("J136").Value = 14499.947
$ws.Range("K136").Value = 7927.5
$ws.Range("L136").Value = 43499.841
$ws.Range("M136").Value = -5377.5
